# Update "projeto de pesquisa" with the new Objetivo / Objetivo Geral /
# Objetivos especificos content.
#
# The document starts out as a single empty paragraph. We turn that
# paragraph into the "Objetivo" heading and then add the remaining
# paragraphs (mixing headings, body text and a couple of blank
# separator paragraphs) right after it, in order, each one inheriting
# the same paragraph formatting (Normal style, bidi=0, left aligned).

$d = $word.ActiveDocument

$paragraphs = @(
    "Objetivo",
    "",
    "Objetivo Geral",
    "O objetivo geral é criar um site para auxiliar imigrantes e brasileiros",
    "",
    "Objetivos específicos",
    "- Auxiliar as pessoas a encontrarem empregos no Brasil",
    "- Normalizar o entendimento referente a legislações e normas brasileiras",
    "- Promover e auxiliar a imigração de mão de obra qualificada no Brasil"
)

# The document already contains exactly one (empty) paragraph; reuse it
# for the first line instead of inserting a new one.
$d.Paragraphs(1).Range.Text = $paragraphs[0]

for ($i = 1; $i -lt $paragraphs.Length; $i++) {
    $d.Paragraphs($i).Range.InsertParagraphAfter()
    if ($paragraphs[$i] -ne "") {
        $d.Paragraphs($i + 1).Range.Text = $paragraphs[$i]
    }
}

# Normal style tweaks that came along with this edit: no hyphenation,
# no paragraph spacing, left-aligned.
$d.Styles("Normal").ParagraphFormat.Hyphenation = $false
$d.Styles("Normal").ParagraphFormat.SpaceBefore = 0
$d.Styles("Normal").ParagraphFormat.SpaceAfter = 0
$d.Styles("Normal").ParagraphFormat.Alignment = 0
